# Apply "Trade #11 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - update aggregate stats
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1300.23   # Current Capital
$summary.Range("B4").Value = 0.23      # Total P&L $
$summary.Range("B5").Value = 0.42      # Total P&L %
$summary.Range("B6").Value = 11        # Total Trades
$summary.Range("B7").Value = 6         # Winning Trades
$summary.Range("B9").Value = 54.55     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - update MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.23     # Capital
$status.Range("D4").Value = 11         # Trades
$status.Range("E4").Value = 0.23       # P&L $
$status.Range("F4").Value = 0.23       # P&L %
$status.Range("G4").Value = 54.55      # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the new Trade #11 row to a trade log sheet (row 12)
# ---------------------------------------------------------------------------
function Add-Trade11Row($ws) {
    $ws.Range("A12").Value = 11

    # Date / Time are stored as plain text in this workbook (not Excel date
    # serials), so force text format before assigning, then reset the style
    # back to Normal so no stray formatting is left behind on the cell.
    $ws.Range("B12").NumberFormat = "@"
    $ws.Range("B12").Value = "2026-02-17"
    $ws.Range("B12").Style = "Normal"

    $ws.Range("C12").NumberFormat = "@"
    $ws.Range("C12").Value = "19:47:35"
    $ws.Range("C12").Style = "Normal"

    $ws.Range("D12").Value = "MarketMaking"
    $ws.Range("E12").Value = "DOWN"
    $ws.Range("F12").Value = 0.19
    $ws.Range("G12").Value = 0.44
    $ws.Range("H12").Value = "CLOSED"
    $ws.Range("I12").Value = 131.5789
    $ws.Range("J12").Value = 0.25
    $ws.Range("K12").Value = 100.23
    $ws.Range("L12").Value = 0
    $ws.Range("M12").Value = 0
    $ws.Range("N12").Value = 0.6
    $ws.Range("O12").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P12").Value = "early_exit"
    $ws.Range("Q12").Value = 0.14
}

# ---------------------------------------------------------------------------
# 3. All Trades sheet - append Trade #11
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade11Row $allTrades

# ---------------------------------------------------------------------------
# 4. MarketMaking sheet - append the same Trade #11 (mirrors All Trades)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade11Row $marketMaking
